# Replace the 100 multiplication-table answers in the document's table with
# their new values (per cell, in document/reading order).
#
# Note: in this runtime, Find.Execute always searches/replaces starting from
# the beginning of $d.Content regardless of which Range invoked .Find, and
# Replace:=1 (normally wdReplaceAll in real Word) is what actually replaces
# only the single first match here, while Replace:=2 (wdReplaceOne) replaces
# every match. So we use Replace:=1 for "replace the next/first occurrence"
# semantics, and we order the individual replacements so that a cell whose
# *new* text would equal another cell's still-unprocessed *old* text is
# always handled before that other cell (otherwise the later Find would hit
# the freshly-written text instead of its own original text).
$d = $word.ActiveDocument

$d.Content.Find.Execute("37×23=851", $true, $false, $false, $false, $false, $true, 1, $false, "85×64=5440", 1) | Out-Null
$d.Content.Find.Execute("75×42=3150", $true, $false, $false, $false, $false, $true, 1, $false, "59×60=3540", 1) | Out-Null
$d.Content.Find.Execute("38×19=722", $true, $false, $false, $false, $false, $true, 1, $false, "81×24=1944", 1) | Out-Null
$d.Content.Find.Execute("96×49=4704", $true, $false, $false, $false, $false, $true, 1, $false, "21×37=777", 1) | Out-Null
$d.Content.Find.Execute("43×50=2150", $true, $false, $false, $false, $false, $true, 1, $false, "63×74=4662", 1) | Out-Null
$d.Content.Find.Execute("86×17=1462", $true, $false, $false, $false, $false, $true, 1, $false, "44×62=2728", 1) | Out-Null
$d.Content.Find.Execute("35×99=3465", $true, $false, $false, $false, $false, $true, 1, $false, "13×42=546", 1) | Out-Null
$d.Content.Find.Execute("46×17=782", $true, $false, $false, $false, $false, $true, 1, $false, "24×74=1776", 1) | Out-Null
$d.Content.Find.Execute("42×69=2898", $true, $false, $false, $false, $false, $true, 1, $false, "46×17=782", 1) | Out-Null
$d.Content.Find.Execute("26×45=1170", $true, $false, $false, $false, $false, $true, 1, $false, "28×41=1148", 1) | Out-Null
$d.Content.Find.Execute("99×97=9603", $true, $false, $false, $false, $false, $true, 1, $false, "57×86=4902", 1) | Out-Null
$d.Content.Find.Execute("63×93=5859", $true, $false, $false, $false, $false, $true, 1, $false, "24×61=1464", 1) | Out-Null
$d.Content.Find.Execute("11×18=198", $true, $false, $false, $false, $false, $true, 1, $false, "54×35=1890", 1) | Out-Null
$d.Content.Find.Execute("70×63=4410", $true, $false, $false, $false, $false, $true, 1, $false, "68×73=4964", 1) | Out-Null
$d.Content.Find.Execute("43×62=2666", $true, $false, $false, $false, $false, $true, 1, $false, "97×88=8536", 1) | Out-Null
$d.Content.Find.Execute("78×62=4836", $true, $false, $false, $false, $false, $true, 1, $false, "67×79=5293", 1) | Out-Null
$d.Content.Find.Execute("63×20=1260", $true, $false, $false, $false, $false, $true, 1, $false, "77×66=5082", 1) | Out-Null
$d.Content.Find.Execute("70×89=6230", $true, $false, $false, $false, $false, $true, 1, $false, "87×54=4698", 1) | Out-Null
$d.Content.Find.Execute("25×71=1775", $true, $false, $false, $false, $false, $true, 1, $false, "79×30=2370", 1) | Out-Null
$d.Content.Find.Execute("27×67=1809", $true, $false, $false, $false, $false, $true, 1, $false, "85×91=7735", 1) | Out-Null
$d.Content.Find.Execute("52×71=3692", $true, $false, $false, $false, $false, $true, 1, $false, "78×63=4914", 1) | Out-Null
$d.Content.Find.Execute("63×53=3339", $true, $false, $false, $false, $false, $true, 1, $false, "75×82=6150", 1) | Out-Null
$d.Content.Find.Execute("60×81=4860", $true, $false, $false, $false, $false, $true, 1, $false, "47×26=1222", 1) | Out-Null
$d.Content.Find.Execute("20×35=700", $true, $false, $false, $false, $false, $true, 1, $false, "26×21=546", 1) | Out-Null
$d.Content.Find.Execute("75×32=2400", $true, $false, $false, $false, $false, $true, 1, $false, "16×57=912", 1) | Out-Null
$d.Content.Find.Execute("23×10=230", $true, $false, $false, $false, $false, $true, 1, $false, "49×16=784", 1) | Out-Null
$d.Content.Find.Execute("13×87=1131", $true, $false, $false, $false, $false, $true, 1, $false, "97×57=5529", 1) | Out-Null
$d.Content.Find.Execute("93×72=6696", $true, $false, $false, $false, $false, $true, 1, $false, "85×60=5100", 1) | Out-Null
$d.Content.Find.Execute("46×22=1012", $true, $false, $false, $false, $false, $true, 1, $false, "43×79=3397", 1) | Out-Null
$d.Content.Find.Execute("85×25=2125", $true, $false, $false, $false, $false, $true, 1, $false, "28×34=952", 1) | Out-Null
$d.Content.Find.Execute("95×82=7790", $true, $false, $false, $false, $false, $true, 1, $false, "35×46=1610", 1) | Out-Null
$d.Content.Find.Execute("100×56=5600", $true, $false, $false, $false, $false, $true, 1, $false, "56×12=672", 1) | Out-Null
$d.Content.Find.Execute("57×25=1425", $true, $false, $false, $false, $false, $true, 1, $false, "73×89=6497", 1) | Out-Null
$d.Content.Find.Execute("72×96=6912", $true, $false, $false, $false, $false, $true, 1, $false, "78×52=4056", 1) | Out-Null
$d.Content.Find.Execute("18×31=558", $true, $false, $false, $false, $false, $true, 1, $false, "65×15=975", 1) | Out-Null
$d.Content.Find.Execute("19×92=1748", $true, $false, $false, $false, $false, $true, 1, $false, "55×31=1705", 1) | Out-Null
$d.Content.Find.Execute("27×99=2673", $true, $false, $false, $false, $false, $true, 1, $false, "18×63=1134", 1) | Out-Null
$d.Content.Find.Execute("53×40=2120", $true, $false, $false, $false, $false, $true, 1, $false, "25×15=375", 1) | Out-Null
$d.Content.Find.Execute("79×53=4187", $true, $false, $false, $false, $false, $true, 1, $false, "96×97=9312", 1) | Out-Null
$d.Content.Find.Execute("23×94=2162", $true, $false, $false, $false, $false, $true, 1, $false, "21×33=693", 1) | Out-Null
$d.Content.Find.Execute("80×52=4160", $true, $false, $false, $false, $false, $true, 1, $false, "66×28=1848", 1) | Out-Null
$d.Content.Find.Execute("86×13=1118", $true, $false, $false, $false, $false, $true, 1, $false, "89×10=890", 1) | Out-Null
$d.Content.Find.Execute("49×52=2548", $true, $false, $false, $false, $false, $true, 1, $false, "81×73=5913", 1) | Out-Null
$d.Content.Find.Execute("84×72=6048", $true, $false, $false, $false, $false, $true, 1, $false, "34×87=2958", 1) | Out-Null
$d.Content.Find.Execute("10×10=100", $true, $false, $false, $false, $false, $true, 1, $false, "73×98=7154", 1) | Out-Null
$d.Content.Find.Execute("88×45=3960", $true, $false, $false, $false, $false, $true, 1, $false, "81×21=1701", 1) | Out-Null
$d.Content.Find.Execute("33×87=2871", $true, $false, $false, $false, $false, $true, 1, $false, "20×60=1200", 1) | Out-Null
$d.Content.Find.Execute("53×10=530", $true, $false, $false, $false, $false, $true, 1, $false, "66×22=1452", 1) | Out-Null
$d.Content.Find.Execute("44×79=3476", $true, $false, $false, $false, $false, $true, 1, $false, "15×75=1125", 1) | Out-Null
$d.Content.Find.Execute("38×58=2204", $true, $false, $false, $false, $false, $true, 1, $false, "73×58=4234", 1) | Out-Null
$d.Content.Find.Execute("64×73=4672", $true, $false, $false, $false, $false, $true, 1, $false, "40×99=3960", 1) | Out-Null
$d.Content.Find.Execute("17×25=425", $true, $false, $false, $false, $false, $true, 1, $false, "70×55=3850", 1) | Out-Null
$d.Content.Find.Execute("78×78=6084", $true, $false, $false, $false, $false, $true, 1, $false, "52×81=4212", 1) | Out-Null
$d.Content.Find.Execute("98×46=4508", $true, $false, $false, $false, $false, $true, 1, $false, "76×88=6688", 1) | Out-Null
$d.Content.Find.Execute("36×53=1908", $true, $false, $false, $false, $false, $true, 1, $false, "27×76=2052", 1) | Out-Null
$d.Content.Find.Execute("76×12=912", $true, $false, $false, $false, $false, $true, 1, $false, "74×42=3108", 1) | Out-Null
$d.Content.Find.Execute("93×39=3627", $true, $false, $false, $false, $false, $true, 1, $false, "35×44=1540", 1) | Out-Null
$d.Content.Find.Execute("54×41=2214", $true, $false, $false, $false, $false, $true, 1, $false, "42×88=3696", 1) | Out-Null
$d.Content.Find.Execute("76×90=6840", $true, $false, $false, $false, $false, $true, 1, $false, "81×70=5670", 1) | Out-Null
$d.Content.Find.Execute("92×27=2484", $true, $false, $false, $false, $false, $true, 1, $false, "61×14=854", 1) | Out-Null
$d.Content.Find.Execute("76×100=7600", $true, $false, $false, $false, $false, $true, 1, $false, "26×89=2314", 1) | Out-Null
$d.Content.Find.Execute("68×67=4556", $true, $false, $false, $false, $false, $true, 1, $false, "82×35=2870", 1) | Out-Null
$d.Content.Find.Execute("86×39=3354", $true, $false, $false, $false, $false, $true, 1, $false, "35×53=1855", 1) | Out-Null
$d.Content.Find.Execute("63×87=5481", $true, $false, $false, $false, $false, $true, 1, $false, "22×41=902", 1) | Out-Null
$d.Content.Find.Execute("72×74=5328", $true, $false, $false, $false, $false, $true, 1, $false, "13×19=247", 1) | Out-Null
$d.Content.Find.Execute("86×63=5418", $true, $false, $false, $false, $false, $true, 1, $false, "57×92=5244", 1) | Out-Null
$d.Content.Find.Execute("52×55=2860", $true, $false, $false, $false, $false, $true, 1, $false, "42×87=3654", 1) | Out-Null
$d.Content.Find.Execute("77×42=3234", $true, $false, $false, $false, $false, $true, 1, $false, "71×87=6177", 1) | Out-Null
$d.Content.Find.Execute("81×72=5832", $true, $false, $false, $false, $false, $true, 1, $false, "76×13=988", 1) | Out-Null
$d.Content.Find.Execute("17×20=340", $true, $false, $false, $false, $false, $true, 1, $false, "72×90=6480", 1) | Out-Null
$d.Content.Find.Execute("95×18=1710", $true, $false, $false, $false, $false, $true, 1, $false, "90×15=1350", 1) | Out-Null
$d.Content.Find.Execute("57×73=4161", $true, $false, $false, $false, $false, $true, 1, $false, "39×100=3900", 1) | Out-Null
$d.Content.Find.Execute("58×90=5220", $true, $false, $false, $false, $false, $true, 1, $false, "68×91=6188", 1) | Out-Null
$d.Content.Find.Execute("59×86=5074", $true, $false, $false, $false, $false, $true, 1, $false, "70×97=6790", 1) | Out-Null
$d.Content.Find.Execute("82×80=6560", $true, $false, $false, $false, $false, $true, 1, $false, "77×67=5159", 1) | Out-Null
$d.Content.Find.Execute("22×49=1078", $true, $false, $false, $false, $false, $true, 1, $false, "71×73=5183", 1) | Out-Null
$d.Content.Find.Execute("94×61=5734", $true, $false, $false, $false, $false, $true, 1, $false, "60×18=1080", 1) | Out-Null
$d.Content.Find.Execute("66×31=2046", $true, $false, $false, $false, $false, $true, 1, $false, "65×52=3380", 1) | Out-Null
$d.Content.Find.Execute("55×29=1595", $true, $false, $false, $false, $false, $true, 1, $false, "42×49=2058", 1) | Out-Null
$d.Content.Find.Execute("43×78=3354", $true, $false, $false, $false, $false, $true, 1, $false, "48×69=3312", 1) | Out-Null
$d.Content.Find.Execute("85×92=7820", $true, $false, $false, $false, $false, $true, 1, $false, "41×21=861", 1) | Out-Null
$d.Content.Find.Execute("29×97=2813", $true, $false, $false, $false, $false, $true, 1, $false, "67×95=6365", 1) | Out-Null
$d.Content.Find.Execute("90×36=3240", $true, $false, $false, $false, $false, $true, 1, $false, "27×24=648", 1) | Out-Null
$d.Content.Find.Execute("69×58=4002", $true, $false, $false, $false, $false, $true, 1, $false, "87×85=7395", 1) | Out-Null
$d.Content.Find.Execute("12×35=420", $true, $false, $false, $false, $false, $true, 1, $false, "36×10=360", 1) | Out-Null
$d.Content.Find.Execute("32×34=1088", $true, $false, $false, $false, $false, $true, 1, $false, "56×15=840", 1) | Out-Null
$d.Content.Find.Execute("14×72=1008", $true, $false, $false, $false, $false, $true, 1, $false, "86×97=8342", 1) | Out-Null
$d.Content.Find.Execute("76×48=3648", $true, $false, $false, $false, $false, $true, 1, $false, "74×15=1110", 1) | Out-Null
$d.Content.Find.Execute("82×23=1886", $true, $false, $false, $false, $false, $true, 1, $false, "47×35=1645", 1) | Out-Null
$d.Content.Find.Execute("82×42=3444", $true, $false, $false, $false, $false, $true, 1, $false, "63×79=4977", 1) | Out-Null
$d.Content.Find.Execute("60×93=5580", $true, $false, $false, $false, $false, $true, 1, $false, "57×49=2793", 1) | Out-Null
$d.Content.Find.Execute("89×67=5963", $true, $false, $false, $false, $false, $true, 1, $false, "22×83=1826", 1) | Out-Null
$d.Content.Find.Execute("83×36=2988", $true, $false, $false, $false, $false, $true, 1, $false, "29×72=2088", 1) | Out-Null
$d.Content.Find.Execute("63×98=6174", $true, $false, $false, $false, $false, $true, 1, $false, "35×66=2310", 1) | Out-Null
$d.Content.Find.Execute("48×75=3600", $true, $false, $false, $false, $false, $true, 1, $false, "17×28=476", 1) | Out-Null
$d.Content.Find.Execute("54×22=1188", $true, $false, $false, $false, $false, $true, 1, $false, "43×52=2236", 1) | Out-Null
$d.Content.Find.Execute("63×64=4032", $true, $false, $false, $false, $false, $true, 1, $false, "23×21=483", 1) | Out-Null
$d.Content.Find.Execute("96×82=7872", $true, $false, $false, $false, $false, $true, 1, $false, "22×45=990", 1) | Out-Null
$d.Content.Find.Execute("36×52=1872", $true, $false, $false, $false, $false, $true, 1, $false, "25×24=600", 1) | Out-Null
$d.Content.Find.Execute("54×34=1836", $true, $false, $false, $false, $false, $true, 1, $false, "73×16=1168", 1) | Out-Null
